# "add upstream sites missed"
# The source workbook had a stray row 46 containing a lone backtick ("`")
# placeholder in column H. That row is removed, and four upstream
# confirmation sites that were missed from the original survey pass are
# appended as new rows (42-45) right after the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray placeholder row (H46 = "`").
$ws.Rows.Item(46).Delete()

# Row 42: reference_number 58 - Gibson Creek / 195288_us
$ws.Cells.Item(42, 1).Value = 58
$ws.Cells.Item(42, 2).Value = "Gibson Creek"
$ws.Cells.Item(42, 3).Value = "195288_us"
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 9).Value = "high"

# Row 43: reference_number 59 - Barren Creek / 197665_us
$ws.Cells.Item(43, 1).Value = 59
$ws.Cells.Item(43, 2).Value = "Barren Creek"
$ws.Cells.Item(43, 3).Value = "197665_us"
$ws.Cells.Item(43, 4).Value = 100
$ws.Cells.Item(43, 9).Value = "moderate"

# Row 44: reference_number 60 - Moan Creek / 197667_us
$ws.Cells.Item(44, 1).Value = 60
$ws.Cells.Item(44, 2).Value = "Moan Creek"
$ws.Cells.Item(44, 3).Value = "197667_us"
$ws.Cells.Item(44, 4).Value = 100
$ws.Cells.Item(44, 9).Value = "high"

# Row 45: reference_number 61 - Coffin Creek / 197668_us
$ws.Cells.Item(45, 1).Value = 61
$ws.Cells.Item(45, 2).Value = "Coffin Creek"
$ws.Cells.Item(45, 3).Value = "197668_us"
$ws.Cells.Item(45, 4).Value = 40
$ws.Cells.Item(45, 9).Value = "high"

# Restore the view: scrolled down a few rows with the last new entry selected.
$ws.Range("G45").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "habitat_confirmations_priorities: added 4 missed upstream sites"
